$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.979.62"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "2.342.97"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "474.23"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.29"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +18.81%  "
$ws.Range("D9").Value = "2.341.73"
$ws.Range("E9").Value = "  -5.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0965"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.41"
$ws.Range("E11").Value = "  -6.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("E12").Value = "  -2.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.124"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "2.749.00"
$ws.Range("E14").Value = "  -5.49%  "
$ws.Range("D15").Value = "55.011.01"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.95"
$ws.Range("E16").Value = "  -4.89%  "
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("D18").Value = "2.346.17"
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.53"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.55"
$ws.Range("E21").Value = "  -4.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.59"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.74"
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.392"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.153"
$ws.Range("E27").Value = "  -5.88%  "
$ws.Range("D28").Value = "2.452.59"
$ws.Range("E28").Value = "  -5.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.10"
$ws.Range("E29").Value = "  -6.67%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "0.0₃0747"
$ws.Range("E31").Value = "  -5.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "146.25"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.06"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.57"
$ws.Range("E36").Value = "  -3.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.08"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.807"
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.59"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0984"
$ws.Range("E41").Value = "  +6.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.32"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.574"
$ws.Range("E44").Value = "  -5.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0517"
$ws.Range("E45").Value = "  -6.92%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "249.78"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.31"
$ws.Range("E49").Value = "  -10.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.62"
$ws.Range("E50").Value = "  -4.98%  "
$ws.Range("D51").Value = "1.770.69"
$ws.Range("E51").Value = "  -5.28%  "
